# "Load Data To DashBoard"
# 1. Rename the existing (only) sheet to "Task Manager".
# 2. Add a new sheet "Phân Chia Công Việc" (work-allocation dashboard) after it.
# 3. Add a small PERMISSION-lookup table to the "Task Manager" sheet (rows 12-13).
# 4. Populate the new "Phân Chia Công Việc" sheet with the team/task-assignment table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Trang_tính1" -> "Task Manager"
# ---------------------------------------------------------------------------
$wsTask = $wb.ActiveSheet
$wsTask.Name = "Task Manager"

# ---------------------------------------------------------------------------
# 2. Add new worksheet "Phân Chia Công Việc" right after "Task Manager"
# ---------------------------------------------------------------------------
$wsWork = $wb.Worksheets.Add($null, $wsTask)
$wsWork.Name = "Phân Chia Công Việc"

# ---------------------------------------------------------------------------
# 3. "Task Manager" sheet: new mini PERMISSION table in rows 12-13 (B:G)
#    (mirrors the IDPermission / IDUser / IDProject / Role / Object / Privilege
#     columns already used by the PERMISSION block higher up on the sheet)
# ---------------------------------------------------------------------------
$wsTask.Range("C11").Clear()

$wsTask.Range("B12:G12").HorizontalAlignment = -4131
$wsTask.Range("B12:G12").VerticalAlignment = -4108
$wsTask.Range("B12:G12").IndentLevel = 3

$wsTask.Range("B12").Value = "IDPermission"
$wsTask.Range("C12").Value = "IDUser"
$wsTask.Range("D12").Value = "IDProject"
$wsTask.Range("E12").Value = "Role"
$wsTask.Range("F12").Value = "Object"
$wsTask.Range("G12").Value = "Privilege"

$wsTask.Range("B13:G13").VerticalAlignment = -4108
$wsTask.Range("B13:G13").WrapText = $true

$wsTask.Range("B13").Value = 1
$wsTask.Range("C13").Value = 1
$wsTask.Range("D13").Value = 1
$wsTask.Range("E13").Value = "Design"
$wsTask.Range("F13").Value = "Task"
$wsTask.Range("G13").Value = "ADD, MOVE"

$wsTask.Range("F16").Select()

# ---------------------------------------------------------------------------
# 4. "Phân Chia Công Việc" sheet content
# ---------------------------------------------------------------------------
$wsWork.Columns.Item(1).ColumnWidth = 24.27
$wsWork.Columns.Item(2).ColumnWidth = 20.63
$wsWork.Columns.Item(3).ColumnWidth = 59.54
$wsWork.Columns.Item(4).ColumnWidth = 20.63

# Header row
$headerRange = $wsWork.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.RowHeight = 45

$wsWork.Range("A1").Value = "HỌ & TÊN"
$wsWork.Range("B1").Value = "MSSV"
$wsWork.Range("C1").Value = "CÔNG VIỆC"
$wsWork.Range("D1").Value = "ĐÁNH GIÁ"

# Data rows 2-6: columns A (name) + B (student id)
$abRange = $wsWork.Range("A2:B6")
$abRange.HorizontalAlignment = -4131
$abRange.VerticalAlignment = -4108
$abRange.WrapText = $true
$abRange.IndentLevel = 3

# Column C (task description)
$cRange = $wsWork.Range("C2:C6")
$cRange.HorizontalAlignment = -4131
$cRange.VerticalAlignment = -4108
$cRange.WrapText = $true
$cRange.IndentLevel = 2

# Column D (evaluation / progress score)
$dRange = $wsWork.Range("D2:D6")
$dRange.VerticalAlignment = -4108
$dRange.WrapText = $true

$wsWork.Range("A2").Value = "Lầu Hoàng Nguyên"
$wsWork.Range("B2").Value = 2180607802
$wsWork.Range("C2").Value = "Task Progress"

$wsWork.Range("A3").Value = "Nguyễn Thành Nhân"
$wsWork.Range("B3").Value = 2180607824
$wsWork.Range("C3").Value = "Đăng Nhập, Đăng Ký, Xác thực và ủy quyền người dùng bằng OAuth 2.0 đăng nhập bằng Google"
$wsWork.Range("D3").Value = 0.85
$wsWork.Range("D3").NumberFormat = "0%"
$wsWork.Range("D3").HorizontalAlignment = -4108

$wsWork.Range("A4").Value = "Nguyễn Văn Phong"
$wsWork.Range("B4").Value = 2180607874
$wsWork.Range("C4").Value = "Sửa Database"

$wsWork.Range("A5").Value = "Võ Đình Thiên Phú"
$wsWork.Range("B5").Value = 2180607177
$wsWork.Range("C5").Value = "Báo Cáo"

$wsWork.Range("A6").Value = "Phạm Văn Phước"
$wsWork.Range("B6").Value = 2180609103
$wsWork.Range("C6").Value = "Chức Năng Chuyển Đổi Status"

# Trailing empty-but-formatted rows (7-10) to match the source table padding
$wsWork.Range("A7:D10").VerticalAlignment = -4108
$wsWork.Range("A7:D10").WrapText = $true

for ($r = 1; $r -le 16; $r++) {
    $wsWork.Rows.Item($r).RowHeight = 30
}
$wsWork.Rows.Item(1).RowHeight = 45
$wsWork.Rows.Item(2).RowHeight = 45
$wsWork.Rows.Item(3).RowHeight = 45
$wsWork.Rows.Item(4).RowHeight = 45
$wsWork.Rows.Item(5).RowHeight = 45
$wsWork.Rows.Item(6).RowHeight = 45

$wsWork.Range("C11").Select()

$wsTask.Activate()
